$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81 (shifts existing rows 81-129 down to 82-130,
# and copies formatting - e.g. the date style on column D - from the row above).
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new daily price record.
$ws.Range("A81").Value = 5
$ws.Range("B81").Value = "Macroferia Regional de Talca"
$ws.Range("C81").Value = "Maule"
$ws.Range("D81").Value = "2022-06-17"
$ws.Range("E81").Value = 7
$ws.Range("F81").Value = 100112001
$ws.Range("G81").Value = "Berenjena"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 300
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = 8000
$ws.Range("N81").Value = "$/caja 50 unidades"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 160
$ws.Range("Q81").Value = 50
$ws.Range("R81").Value = "Hortaliza"
